$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = "Giovanni Simoncelli"
$ws.Range("B52").Value = "Stefano  Tita | Clitoriders"
$ws.Range("C52").Value = "Daniele Feltrinelli | Rita Levi’s"
$ws.Range("D52").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E52").Value = "Andrea  Pedrotti | IMONTAGNA"
$ws.Range("F52").Value = "Alessandro Galvagni | Hellas Lazio"
